$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# 1. Append " Class 6R drivers licence." to the end of the personal
#    statement paragraph ("...business networking.").
# -------------------------------------------------------------------------
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$found = $findRange.Find.Execute("networking.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$findRange.Collapse(0)
$insertStart = $findRange.Start

# Insert the sentence plus one throw-away padding character. The padding
# keeps the bookmark-target position (computed below) from sitting exactly
# on the paragraph-mark boundary while we create the bookmark; it is
# deleted again immediately afterwards.
$findRange.InsertAfter(" Class 6R drivers licence.X")

# The inserted text lands in the same run as the preceding sentence
# (identical character formatting), so split it into its own run the way
# Word does when a run boundary is forced - toggle a character property on
# and back off again.
$newRun = $d.Range($insertStart, $insertStart + 26)
$newRun.Bold = 1
$newRun.Bold = 0

# -------------------------------------------------------------------------
# 2. Relocate the "_GoBack" bookmark (last-edit marker) from its old spot
#    next to "LoRa" in the skills list to the end of the text just typed.
#    Re-adding a bookmark under an existing name moves it, so the old
#    bookmark near "LoRa" is automatically replaced.
# -------------------------------------------------------------------------
$bookmarkPos = $insertStart + 26
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Remove the padding character now that the bookmark is safely placed.
$padding = $d.Range($insertStart + 26, $insertStart + 27)
$padding.Delete()
